# Apply the changes described by the diff to the workbook.
$wb = $excel.ActiveWorkbook

# --- Rename the second worksheet ---
$wsInclude = $wb.Worksheets.Item("Include from Health Data Conn")
$wsInclude.Name = "Include from Payer Claim Type"

# --- Update the Metadata sheet values ---
$wsMeta = $wb.Worksheets.Item("Metadata")

# URL
$wsMeta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/ValueSet/claim-type"

# Version
$wsMeta.Range("B3").Value = "8.0.0"

# Date
$wsMeta.Range("B8").Value = "2022-11-10T16:00:46+00:00"

# Publisher
$wsMeta.Range("B9").Value = "LinuxForHealth Team"

# Description
$wsMeta.Range("B11").Value = "Value set for the coverage type under which the claim was paid"

# --- Update the System URI on the renamed "Include from Payer Claim Type" sheet ---
$wsInclude.Range("B4").Value = "http://linuxforhealth.org/fhir/cdm/CodeSystem/payer-claim-type"

$wb.Save()
